$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as plain text matching the original (no style residue),
# used for D-column price cells that could otherwise be parsed as numbers.
function Set-TextCell($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '26.907.64'
$ws.Range('E2').Value = '  -0.14%  '
Set-TextCell $ws 'D3' '1.816.26'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.50%  '
Set-TextCell $ws 'D5' '308.97'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  +0.33%  '
Set-TextCell $ws 'D7' '0.4669'
$ws.Range('E7').Value = '  +0.75%  '
Set-TextCell $ws 'D8' '0.3689'
$ws.Range('E8').Value = '  -1.40%  '
Set-TextCell $ws 'D9' '0.07357'
$ws.Range('E9').Value = '  +0.61%  '
Set-TextCell $ws 'D10' '0.8711'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E11').Value = '  -0.36%  '
Set-TextCell $ws 'D12' '1.790.51'
$ws.Range('E12').Value = '  +2.70%  '
Set-TextCell $ws 'D13' '5.376'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws 'D14' '6.515'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D15' '0.07067'
$ws.Range('E15').Value = '  +0.33%  '
Set-TextCell $ws 'D16' '91.59'
$ws.Range('E16').Value = '  +0.04%  '
Set-TextCell $ws 'D17' '1.003'
$ws.Range('E17').Value = '  +0.57%  '
Set-TextCell $ws 'D18' '0.000008699'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('E20').Value = '  -0.09%  '
Set-TextCell $ws 'D21' '26.938.43'
$ws.Range('E21').Value = '  +0.11%  '
Set-TextCell $ws 'D22' '5.321'
$ws.Range('E22').Value = '  +0.28%  '
Set-TextCell $ws 'D23' '10.61'
$ws.Range('E23').Value = '  -0.86%  '
Set-TextCell $ws 'D24' '2.044.87'
$ws.Range('E24').Value = '  +4.12%  '
Set-TextCell $ws 'D25' '1.896'
$ws.Range('E25').Value = '  -0.67%  '
Set-TextCell $ws 'D26' '150.49'
$ws.Range('E26').Value = '  -0.11%  '
Set-TextCell $ws 'D27' '2.168'
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('E28').Value = '  +0.01%  '
Set-TextCell $ws 'D29' '5.334'
$ws.Range('E29').Value = '  +1.38%  '
Set-TextCell $ws 'D30' '115.87'
$ws.Range('E30').Value = '  +0.92%  '
Set-TextCell $ws 'D31' '0.08947'
$ws.Range('E31').Value = '  +0.51%  '
Set-TextCell $ws 'D32' '0.7686'
$ws.Range('E32').Value = '  -0.02%  '
Set-TextCell $ws 'D33' '1.163'
$ws.Range('E33').Value = '  -0.89%  '
Set-TextCell $ws 'D34' '4.504'
$ws.Range('E34').Value = '  +0.77%  '
Set-TextCell $ws 'D35' '2.905'
$ws.Range('E35').Value = '  +0.62%  '
Set-TextCell $ws 'D36' '1.001'
$ws.Range('E36').Value = '  +0.44%  '
Set-TextCell $ws 'D37' '1.086'
$ws.Range('E37').Value = '  -3.09%  '
Set-TextCell $ws 'D38' '0.01964'
$ws.Range('E38').Value = '  +0.64%  '
Set-TextCell $ws 'D39' '0.05290'
$ws.Range('E39').Value = '  +1.01%  '
Set-TextCell $ws 'D40' '2.931'
$ws.Range('E40').Value = '  +0.89%  '
Set-TextCell $ws 'D41' '7.260'
$ws.Range('E41').Value = '  +0.87%  '
Set-TextCell $ws 'D42' '0.5318'
$ws.Range('E42').Value = '  +1.23%  '
Set-TextCell $ws 'D43' '2.350'
$ws.Range('E43').Value = '  -3.78%  '
Set-TextCell $ws 'D44' '0.1658'
$ws.Range('E44').Value = '  -0.15%  '
Set-TextCell $ws 'D45' '8.440'
$ws.Range('E45').Value = '  -1.64%  '
Set-TextCell $ws 'D46' '0.4924'
$ws.Range('E46').Value = '  -2.61%  '
Set-TextCell $ws 'D47' '10.42'
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell $ws 'D48' '1.001'
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D49' '1.671'
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('E50').Value = '  -0.57%  '
Set-TextCell $ws 'D51' '0.06300'
$ws.Range('E51').Value = '  -0.18%  '
